# Applies the diff to the speaker notes of slide 14 ("Evaluation Metrics")
# and slide 21 ("MOT dataset screenshots").
#
# NB: PowerPoint speaker-notes TextRange objects on this host only support
# whole-body assignment (TextRange.Text = "..."), so each edit below
# rebuilds the full notes body, joining paragraphs with an LF so the host
# splits them into separate <a:p> paragraphs.

$p = $ppt.ActivePresentation
$nl = [char]10

# --- Slide 14: add explanatory sentence about green/orange coloring ----
$s14 = $p.Slides.Item(14)
$notes14 = $s14.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes14.Text = "Green ones are those where higher is better while orange ones have lower is better;" + $nl

# --- Slide 21: add notes about running the code and the evaluation kit -
$s21 = $p.Slides.Item(21)
$notes21 = $s21.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes21.Text = "These are some screenshots from the MOT dataset with tracked objects shown;" + $nl + "I did run the code on all test sequences from MOT 2015 but it only produces a text file as output that can be read by the MOT evaluation fkit;" + $nl + "This kit, however, provides no tools for visualization;" + $nl + "" + $nl + ""
